# edit.ps1 - applies the "Doc observaciones: descripcion maquina 1" change
# Target changes (content-level; style-id renames / gridCol resizing / page-break
# relocation seen in the source diff are Word-internal recalculation artifacts
# that are not meaningfully reproducible -- or required -- via scripted edits):
#   1. "1 Cod XXXX" -> "1 Cod " + "202110516" (kept as two runs)
#   2. Procesadores / Maquina 1 cell gets the CPU description, jc=both removed
#   3. Memoria RAM / Maquina 1 cell gets "8 GB" (jc=both stays)
#   4. Sistema Operativo / Maquina 1 cell gets the OS description (several runs
#      + proofErr markers), jc=both removed

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1. "Estudiante 1 Cod XXXX" paragraph -> split "1 Cod XXXX" into two runs
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(2)
$body1 = '<w:body><w:p><w:pPr><w:spacing w:after="0"/><w:jc w:val="right"/><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-419"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">Estudiante </w:t></w:r>' `
  + '<w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">1 Cod </w:t></w:r>' `
  + '<w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-419"/></w:rPr><w:t>202110516</w:t></w:r>' `
  + '</w:p></w:body>'
$p1.Range.InsertXML($pkgHeader + $body1 + $pkgFooter)

# ---------------------------------------------------------------------------
# Locate the first "Ambientes de pruebas" table (Maquina 1 / Maquina 2 specs)
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# 2. Procesadores / Maquina 1 (row 2, col 2) - add CPU text, drop jc=both
# ---------------------------------------------------------------------------
$cellCpu = $t.Cell(2, 2)
$bodyCpu = '<w:body><w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:noProof w:val="0"/><w:lang w:val="es-419"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:noProof w:val="0"/><w:lang w:val="es-419"/></w:rPr><w:t>AMD A9-9425 RADEON R5, 5 COMPUTE CORES 2C+3G 3.10 GHz</w:t></w:r>' `
  + '</w:p></w:body>'
$cellCpu.Range.InsertXML($pkgHeader + $bodyCpu + $pkgFooter)

# ---------------------------------------------------------------------------
# 3. Memoria RAM / Maquina 1 (row 3, col 2) - add "8 GB", jc=both is kept
# ---------------------------------------------------------------------------
$cellRam = $t.Cell(3, 2)
$bodyRam = '<w:body><w:p><w:pPr><w:jc w:val="both"/><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:noProof w:val="0"/><w:lang w:val="es-419"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:noProof w:val="0"/><w:lang w:val="es-419"/></w:rPr><w:t>8 GB</w:t></w:r>' `
  + '</w:p></w:body>'
$cellRam.Range.InsertXML($pkgHeader + $bodyRam + $pkgFooter)

# ---------------------------------------------------------------------------
# 4. Sistema Operativo / Maquina 1 (row 4, col 2) - add OS text, drop jc=both
# ---------------------------------------------------------------------------
$cellOs = $t.Cell(4, 2)
$rPrOs = '<w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:noProof w:val="0"/><w:lang w:val="es-419"/></w:rPr>'
$bodyOs = '<w:body><w:p><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:noProof w:val="0"/><w:lang w:val="es-419"/></w:rPr></w:pPr>' `
  + '<w:r>' + $rPrOs + '<w:t>Windows</w:t></w:r>' `
  + '<w:r>' + $rPrOs + '<w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:r>' + $rPrOs + '<w:t xml:space="preserve">10 Home </w:t></w:r>' `
  + '<w:proofErr w:type="gramStart"/>' `
  + '<w:r>' + $rPrOs + '<w:t>Single</w:t></w:r>' `
  + '<w:proofErr w:type="gramEnd"/>' `
  + '<w:r>' + $rPrOs + '<w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:r>' + $rPrOs + '<w:t>64-bits</w:t></w:r>' `
  + '</w:p></w:body>'
$cellOs.Range.InsertXML($pkgHeader + $bodyOs + $pkgFooter)

Write-Output "edits applied"
